# Append two new rows (76 and 77) of usage data to the "Data" worksheet,
# mirroring the pattern of the existing rows (column A = timestamp serial,
# columns B:O = numeric counters).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$rows = @(
    @(76, 45806.913958333331, 11, 6, 395, 618, 590, 687, 5536, 687, 58, 6, 682, 30, 5723, 6939),
    @(77, 45807.924907407411, 12, 6, 400, 625, 597, 688, 5567, 688,  2, 2, 684, 30, 5771, 6960)
)

foreach ($row in $rows) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
    $ws.Cells.Item($r, 1).HorizontalAlignment = -4108
    $ws.Cells.Item($r, 1).VerticalAlignment = -4108

    for ($c = 2; $c -le 15; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c]
        $ws.Cells.Item($r, $c).HorizontalAlignment = -4108
        $ws.Cells.Item($r, $c).VerticalAlignment = -4108
    }
}
